$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 and 35 swap (ARBITRUM <-> MXToken) plus new Price/Volume values
$ws.Range("B34").Value = 'MXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.794'
$ws.Range("E34").Value = '  +0.18%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9490'
$ws.Range("E35").Value = '  -2.81%  '

# Price / Volume(1h) updates for the remaining rows
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.119.31'
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.655.46'
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.11'
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5291'
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2612'
$ws.Range("E8").Value = '  -2.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06338'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.42'
$ws.Range("E10").Value = '  -3.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07771'
$ws.Range("E11").Value = '  -0.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.494'
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.654.33'
$ws.Range("E13").Value = '  -0.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5472'
$ws.Range("E14").Value = '  -1.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0₅8163'
$ws.Range("E15").Value = '  -1.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.39'
$ws.Range("E16").Value = '  +0.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.134.66'
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("E18").Value = '  -0.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.561'
$ws.Range("E19").Value = '  -2.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '192.82'
$ws.Range("E20").Value = '  -1.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.09'
$ws.Range("E21").Value = '  -1.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.025'
$ws.Range("E22").Value = '  -1.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.004'
$ws.Range("E23").Value = '  -0.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '141.86'
$ws.Range("E24").Value = '  +1.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1251'
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.269'
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.18'
$ws.Range("E27").Value = '  -0.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.438'
$ws.Range("E28").Value = '  +1.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05934'
$ws.Range("E29").Value = '  -4.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.280'
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.518'
$ws.Range("E31").Value = '  -2.53%  '
$ws.Range("E32").Value = '  -1.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.570'
$ws.Range("E33").Value = '  -4.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.412'
$ws.Range("E36").Value = '  -0.62%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5660'
$ws.Range("E37").Value = '  -2.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01610'
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.807'
$ws.Range("E39").Value = '  -3.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8486'
$ws.Range("E40").Value = '  -1.25%  '
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '102.64'
$ws.Range("E42").Value = '  +2.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.022.88'
$ws.Range("E43").Value = '  -0.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.799.45'
$ws.Range("E44").Value = '  -0.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.18'
$ws.Range("E45").Value = '  -1.34%  '
$ws.Range("E46").Value = '  -0.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4288'
$ws.Range("E47").Value = '  +1.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.477'
$ws.Range("E48").Value = '  -1.09%  '
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.769'
$ws.Range("E50").Value = '  -4.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.09704'
$ws.Range("E51").Value = '  -1.21%  '
